$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Means"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Means")

# New header cells for the "Within 5 miles" / "Within 10 miles" columns
$ws1.Range("F1").Value = "Within 5 miles of HFC production facility"
$ws1.Range("G1").Value = "Within 10 miles of HFC production facility"

# New column F & G data (rows 2-10)
$ws1.Range("F2").Value = 70
$ws1.Range("G2").Value = 64

$ws1.Range("F3").Value = 17
$ws1.Range("G3").Value = 19

$ws1.Range("F4").Value = 12
$ws1.Range("G4").Value = 18

$ws1.Range("F5").Value = 34
$ws1.Range("G5").Value = 27

$ws1.Range("F6").Value = 82
$ws1.Range("G6").Value = 99

$ws1.Range("F7").Value = 5.1
$ws1.Range("G7").Value = 3.5

$ws1.Range("F8").Value = 3.7
$ws1.Range("G8").Value = 3

$ws1.Range("F9").Value = 30
$ws1.Range("G9").Value = 31

$ws1.Range("F10").Value = 0.37
$ws1.Range("G10").Value = 0.37

# Updated values for "Total Cancer Risk (per million)" row (row 9, cols B-E)
$ws1.Range("B9").Value = 26
$ws1.Range("C9").Value = 28
$ws1.Range("D9").Value = 30
$ws1.Range("E9").Value = 30

# Updated values for "Total Respiratory (hazard quotient)" row (row 10, cols B-E)
$ws1.Range("B10").Value = 0.32
$ws1.Range("C10").Value = 0.33
$ws1.Range("D10").Value = 0.35
$ws1.Range("E10").Value = 0.36

# ---------------------------------------------------------------------------
# Sheet 2: "Standard Deviations"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Standard Deviations")

# New header cells for the "Within 5 mile" / "Within 10 mile" SD columns
$ws2.Range("F1").Value = "Within 5 mile of HFC production facility SD"
$ws2.Range("G1").Value = "Within 10 mile of HFC production facility SD"

# New column F & G data (rows 2-10)
$ws2.Range("F2").Value = 14
$ws2.Range("G2").Value = 20

$ws2.Range("F3").Value = 10
$ws2.Range("G3").Value = 15

$ws2.Range("F4").Value = 7.3
$ws2.Range("G4").Value = 10

$ws2.Range("F5").Value = 19
$ws2.Range("G5").Value = 19

$ws2.Range("F6").Value = 25
$ws2.Range("G6").Value = 34

$ws2.Range("F7").Value = 7.1
$ws2.Range("G7").Value = 6

$ws2.Range("F8").Value = 8.7
$ws2.Range("G8").Value = 5.8

$ws2.Range("F9").Value = 0
$ws2.Range("G9").Value = 3.7

$ws2.Range("F10").Value = 0.048
$ws2.Range("G10").Value = 0.044

# Updated values for "Total Cancer Risk (per million)" row (row 9, cols B-E)
$ws2.Range("B9").Value = 8.6
$ws2.Range("C9").Value = 7.5
$ws2.Range("D9").Value = 0
$ws2.Range("E9").Value = 0

# Updated values for "Total Respiratory (hazard quotient)" row (row 10, cols B-E)
$ws2.Range("B10").Value = 0.14
$ws2.Range("C10").Value = 0.076
$ws2.Range("D10").Value = 0.058
$ws2.Range("E10").Value = 0.052
